$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.477.94"
$ws.Range("E2").Value = "'  -2.90%  "
$ws.Range("D3").Value = "'1.798.03"
$ws.Range("E3").Value = "'  -2.36%  "
$ws.Range("E4").Value = "'  +0.35%  "
$ws.Range("D5").Value = "'228.75"
$ws.Range("E5").Value = "'  -1.22%  "
$ws.Range("D6").Value = "'0.611"
$ws.Range("E6").Value = "'  -1.31%  "
$ws.Range("E7").Value = "'  +0.38%  "
$ws.Range("D8").Value = "'38.96"
$ws.Range("E8").Value = "'  -10.69%  "
$ws.Range("D9").Value = "'0.318"
$ws.Range("E9").Value = "'  +2.57%  "
$ws.Range("E10").Value = "'  -3.32%  "
$ws.Range("D11").Value = "'0.0987"
$ws.Range("E11").Value = "'  -2.29%  "
$ws.Range("D12").Value = "'2.060.15"
$ws.Range("B13").Value = "'Chainlink"
$ws.Range("C13").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'11.03"
$ws.Range("E13").Value = "'  -2.13%  "
$ws.Range("B14").Value = "'Polygon"
$ws.Range("C14").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.655"
$ws.Range("E14").Value = "'  -2.62%  "
$ws.Range("B15").Value = "'WrappedEther"
$ws.Range("C15").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "'1.780.72"
$ws.Range("E15").Value = "'  -3.23%  "
$ws.Range("E16").Value = "'  -3.80%  "
$ws.Range("D17").Value = "'34.356.75"
$ws.Range("E17").Value = "'  -3.24%  "
$ws.Range("D18").Value = "'68.85"
$ws.Range("E18").Value = "'  -2.20%  "
$ws.Range("D19").Value = "'0.0₃0776"
$ws.Range("E19").Value = "'  -3.01%  "
$ws.Range("D20").Value = "'239.22"
$ws.Range("E20").Value = "'  -2.12%  "
$ws.Range("D21").Value = "'11.73"
$ws.Range("E21").Value = "'  -2.71%  "
$ws.Range("D22").Value = "'4.65"
$ws.Range("E22").Value = "'  +0.42%  "
$ws.Range("E23").Value = "'  +0.31%  "
$ws.Range("E24").Value = "'  -0.20%  "
$ws.Range("D25").Value = "'172.68"
$ws.Range("E25").Value = "'  +0.75%  "
$ws.Range("D26").Value = "'7.66"
$ws.Range("E26").Value = "'  -4.44%  "
$ws.Range("D27").Value = "'17.11"
$ws.Range("E27").Value = "'  -3.93%  "
$ws.Range("E28").Value = "'  -0.89%  "
$ws.Range("E29").Value = "'  -4.79%  "
$ws.Range("E30").Value = "'  +0.26%  "
$ws.Range("D31").Value = "'3.98"
$ws.Range("E31").Value = "'  +1.22%  "
$ws.Range("E32").Value = "'  -2.07%  "
$ws.Range("D33").Value = "'3.87"
$ws.Range("E33").Value = "'  -5.61%  "
$ws.Range("D34").Value = "'1.22"
$ws.Range("E34").Value = "'  +7.50%  "
$ws.Range("E35").Value = "'  -3.37%  "
$ws.Range("D36").Value = "'0.688"
$ws.Range("E36").Value = "'  -0.41%  "
$ws.Range("D37").Value = "'90.41"
$ws.Range("E37").Value = "'  -4.50%  "
$ws.Range("E38").Value = "'  +4.35%  "
$ws.Range("D39").Value = "'1.315.15"
$ws.Range("E39").Value = "'  -2.42%  "
$ws.Range("E40").Value = "'  -2.84%  "
$ws.Range("D41").Value = "'0.953"
$ws.Range("E41").Value = "'  -5.78%  "
$ws.Range("D42").Value = "'2.41"
$ws.Range("E42").Value = "'  -2.08%  "
$ws.Range("D43").Value = "'14.16"
$ws.Range("E43").Value = "'  -8.16%  "
$ws.Range("D44").Value = "'2.19"
$ws.Range("E44").Value = "'  -10.47%  "
$ws.Range("D45").Value = "'2.71"
$ws.Range("E45").Value = "'  -3.67%  "
$ws.Range("D46").Value = "'6.14"
$ws.Range("E46").Value = "'  -1.90%  "
$ws.Range("E47").Value = "'  -1.18%  "
$ws.Range("D48").Value = "'1.982.85"
$ws.Range("E48").Value = "'  -1.48%  "
$ws.Range("E49").Value = "'  +0.31%  "
$ws.Range("D50").Value = "'0.0658"
$ws.Range("E50").Value = "'  +3.64%  "
$ws.Range("D51").Value = "'97.33"
$ws.Range("E51").Value = "'  -4.96%  "
